$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns A..K
$headers = @("Job_Id", "Job_Title", "Job_Description", "Total_Years_Min_Exp", "Total_Years_Max_Exp", "Work_Mode", "Job_Location", "LinkedIn_Poster", "LinkedIn_Posted", "Resume_received", "Resume_downloaded")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data row (row 2)
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Remote"
$ws.Range("G2").Value = "Bengaluru, Karnataka, India"

# Undo the implicit row auto-height bump caused by the embedded line break
$ws.Rows.Item(2).EntireRow.AutoFit()

